$d = $word.ActiveDocument

# The Title, Author, and Abstract paragraphs each had their text split
# across many single-word runs (with separate space runs in between).
# Collapse each of those paragraphs' runs into a single run containing
# the same text, by re-finding and replacing the full paragraph text.

$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Answers: Introduction to sigma notation", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Answers: Introduction to sigma notation", 2)

$authorRange = $d.Paragraphs(2).Range
$authorRange.Find.Execute("Ifan Howells-Baines, Mark Toner", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Ifan Howells-Baines, Mark Toner", 2)

$abstractRange = $d.Paragraphs(4).Range
$abstractRange.Find.Execute("Answers to questions relating to the guide on introduction to sigma notation.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Answers to questions relating to the guide on introduction to sigma notation.", 2)
